# Update faturamento diario lojas data for BIBI
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Bibi Cell Mundi
$ws.Range("W2").Value = 16452.6
$ws.Range("AG2").Value = 255048.39

# Row 3 - Bibi Cell Vieiralves
$ws.Range("W3").Value = 5538
$ws.Range("AG3").Value = 144110.6

# Row 4 - Bibi Cell Ponta Negra
$ws.Range("V4").Value = 1845
$ws.Range("W4").Value = 2804.7
$ws.Range("AG4").Value = 73566.3

# Row 5 - Bibi Cell Manauara
$ws.Range("V5").Value = 2645.5
$ws.Range("W5").Value = 2351.99
$ws.Range("AG5").Value = 68978.67999999999

# Row 6 - total
$ws.Range("V6").Value = 25994.16
$ws.Range("W6").Value = 27147.29
$ws.Range("AG6").Value = 541703.97
